$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 4000
$ws.Range("I19").Value = 4000
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 4000
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -3825
# Row 70
$ws.Range("H70").Value = 5312
$ws.Range("I70").Value = 4311.222
$ws.Range("J70").Value = 6212.7
$ws.Range("K70").Value = 12933.666
$ws.Range("L70").Value = 18638.1
$ws.Range("M70").Value = -12663.666
$ws.Range("N70").Value = -19178.1
# Row 73
$ws.Range("H73").Value = 5312
$ws.Range("I73").Value = 4311.222
$ws.Range("J73").Value = 6212.7
$ws.Range("K73").Value = 12933.666
$ws.Range("L73").Value = 18638.1
$ws.Range("M73").Value = -11997.666
$ws.Range("N73").Value = -20510.1
# Row 132
$ws.Range("H132").Value = 11775.7
$ws.Range("I132").Value = 12206.333
$ws.Range("J132").Value = 7900
$ws.Range("K132").Value = 36618.999
$ws.Range("L132").Value = 23700
$ws.Range("M132").Value = -34088.999
$ws.Range("N132").Value = -28760

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2094.875
$ws.Range("I2").Value = 2280.4285
$ws.Range("J2").Value = 796
$ws.Range("K2").Value = 2280.4285
$ws.Range("L2").Value = 796
$ws.Range("M2").Value = -2167.4285
$ws.Range("N2").Value = -1022
# Row 32
$ws.Range("H32").Value = 3572174.2
$ws.Range("I32").Value = 609.125
$ws.Range("J32").Value = 25001566
$ws.Range("K32").Value = 609.125
$ws.Range("L32").Value = 25001566
$ws.Range("M32").Value = -322.125
$ws.Range("N32").Value = -25002140
# Row 61
$ws.Range("H61").Value = 2218.3333
$ws.Range("I61").Value = 1562.6
$ws.Range("J61").Value = 5497
$ws.Range("K61").Value = 1562.6
$ws.Range("L61").Value = 5497
$ws.Range("M61").Value = -1350.6
$ws.Range("N61").Value = -5921
# Row 112
$ws.Range("H112").Value = 2387
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 2387
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 2387
$ws.Range("N112").Value = -5341
# Row 116
$ws.Range("H116").Value = 2094.875
$ws.Range("I116").Value = 2280.4285
$ws.Range("J116").Value = 796
$ws.Range("K116").Value = 2280.4285
$ws.Range("L116").Value = 796
$ws.Range("M116").Value = 13.57150000000001
$ws.Range("N116").Value = -5384
# Row 124
$ws.Range("H124").Value = 25724.25
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 25724.25
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 25724.25
$ws.Range("N124").Value = -35544.25
# Row 125
$ws.Range("H125").Value = 71500
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 71500
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 71500
$ws.Range("N125").Value = -81340
# Row 132
$ws.Range("H132").Value = 2110
$ws.Range("I132").Value = 2110
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6330
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3800
# Row 135
$ws.Range("H135").Value = 85249.25
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 85249.25
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 85249.25
$ws.Range("N135").Value = -95389.25
# Row 136
$ws.Range("H136").Value = 2218.3333
$ws.Range("I136").Value = 1562.6
$ws.Range("J136").Value = 5497
$ws.Range("K136").Value = 4687.799999999999
$ws.Range("L136").Value = 16491
$ws.Range("M136").Value = -2137.799999999999
$ws.Range("N136").Value = -21591

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2094.875
$ws.Range("I3").Value = 2280.4285
$ws.Range("J3").Value = 796
$ws.Range("K3").Value = 2280.4285
$ws.Range("L3").Value = 796
$ws.Range("M3").Value = -2166.4285
$ws.Range("N3").Value = -1024
# Row 20
$ws.Range("H20").Value = 4437.857
$ws.Range("I20").Value = 4413.2
$ws.Range("J20").Value = 4499.5
$ws.Range("K20").Value = 4413.2
$ws.Range("L20").Value = 4499.5
$ws.Range("M20").Value = -4166.2
$ws.Range("N20").Value = -4993.5
# Row 86
$ws.Range("H86").Value = 4665.533
$ws.Range("I86").Value = 2098.3
$ws.Range("J86").Value = 9800
$ws.Range("K86").Value = 2098.3
$ws.Range("L86").Value = 9800
$ws.Range("M86").Value = -975.3000000000002
$ws.Range("N86").Value = -12046
# Row 89
$ws.Range("H89").Value = 4665.533
$ws.Range("I89").Value = 2098.3
$ws.Range("J89").Value = 9800
$ws.Range("K89").Value = 10491.5
$ws.Range("L89").Value = 49000
$ws.Range("M89").Value = -4875.5
$ws.Range("N89").Value = -60232
# Row 94
$ws.Range("H94").Value = 179.91667
$ws.Range("I94").Value = 143
$ws.Range("J94").Value = 364.5
$ws.Range("K94").Value = 143
$ws.Range("L94").Value = 364.5
$ws.Range("M94").Value = 308
$ws.Range("N94").Value = -1266.5
# Row 99
$ws.Range("H99").Value = 500001000
$ws.Range("I99").Value = 500001000
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 500001000
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -499999502
# Row 100
$ws.Range("H100").Value = 14900
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 14900
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 14900
$ws.Range("N100").Value = -17064
# Row 103
$ws.Range("H103").Value = 11938
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 11938
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 11938
$ws.Range("N103").Value = -14282

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 2801.8823
$ws.Range("I58").Value = 2017.2858
$ws.Range("J58").Value = 6463.3335
$ws.Range("K58").Value = 2017.2858
$ws.Range("L58").Value = 6463.3335
$ws.Range("M58").Value = -1814.2858
$ws.Range("N58").Value = -6869.3335
# Row 99
$ws.Range("H99").Value = 2139.4
$ws.Range("I99").Value = 2061.875
$ws.Range("J99").Value = 2449.5
$ws.Range("K99").Value = 2061.875
$ws.Range("L99").Value = 2449.5
$ws.Range("M99").Value = -563.875
$ws.Range("N99").Value = -5445.5
# Row 122
$ws.Range("H122").Value = 1750
$ws.Range("I122").Value = 1750
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5250
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2800
# Row 126
$ws.Range("H126").Value = 2139.4
$ws.Range("I126").Value = 2061.875
$ws.Range("J126").Value = 2449.5
$ws.Range("K126").Value = 6185.625
$ws.Range("L126").Value = 7348.5
$ws.Range("M126").Value = -3715.625
$ws.Range("N126").Value = -12288.5
# Row 136
$ws.Range("H136").Value = 2801.8823
$ws.Range("I136").Value = 2017.2858
$ws.Range("J136").Value = 6463.3335
$ws.Range("K136").Value = 6051.857400000001
$ws.Range("L136").Value = 19390.0005
$ws.Range("M136").Value = -3501.857400000001
$ws.Range("N136").Value = -24490.0005

$ws = $wb.Worksheets.Item("CUL")
# Row 35
$ws.Range("H35").Value = 159
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 159
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 477
$ws.Range("N35").Value = -1053
# Row 80
$ws.Range("H80").Value = 4327.2085
$ws.Range("I80").Value = 4066.611
$ws.Range("J80").Value = 5109
$ws.Range("K80").Value = 12199.833
$ws.Range("L80").Value = 15327
$ws.Range("M80").Value = -11263.833
$ws.Range("N80").Value = -17199
# Row 82
$ws.Range("H82").Value = 7000
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 7000
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 21000
$ws.Range("N82").Value = -21812
# Row 83
$ws.Range("H83").Value = 4327.2085
$ws.Range("I83").Value = 4066.611
$ws.Range("J83").Value = 5109
$ws.Range("K83").Value = 36599.499
$ws.Range("L83").Value = 45981
$ws.Range("M83").Value = -31919.499
$ws.Range("N83").Value = -55341
# Row 85
$ws.Range("H85").Value = 7000
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 7000
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 21000
$ws.Range("N85").Value = -23808

$ws = $wb.Worksheets.Item("GSM")
# Row 57
$ws.Range("H57").Value = 70706
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 70706
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 70706
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -72346
# Row 80
$ws.Range("H80").Value = 1930.1538
$ws.Range("I80").Value = 1756.8572
$ws.Range("J80").Value = 2132.3333
$ws.Range("K80").Value = 1756.8572
$ws.Range("L80").Value = 2132.3333
$ws.Range("M80").Value = -758.8571999999999
$ws.Range("N80").Value = -4128.3333
# Row 83
$ws.Range("H83").Value = 1930.1538
$ws.Range("I83").Value = 1756.8572
$ws.Range("J83").Value = 2132.3333
$ws.Range("K83").Value = 8784.286
$ws.Range("L83").Value = 10661.6665
$ws.Range("M83").Value = -3792.286
$ws.Range("N83").Value = -20645.6665
# Row 92
$ws.Range("H92").Value = 9999.25
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 9999.25
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 9999.25
$ws.Range("N92").Value = -13743.25
# Row 132
$ws.Range("H132").Value = 4622.5557
$ws.Range("I132").Value = 4806.2354
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 14418.7062
$ws.Range("L132").Value = 4500
$ws.Range("M132").Value = -11888.7062
$ws.Range("N132").Value = -9560

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1142.8572
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 1200
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 1200
$ws.Range("M22").Value = -705
$ws.Range("N22").Value = -1790
# Row 27
$ws.Range("H27").Value = 1142.8572
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 1200
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 1200
$ws.Range("M27").Value = -893
$ws.Range("N27").Value = -1414
# Row 46
$ws.Range("H46").Value = 7015.8335
$ws.Range("I46").Value = 3200
$ws.Range("J46").Value = 7779
$ws.Range("K46").Value = 3200
$ws.Range("L46").Value = 7779
$ws.Range("M46").Value = -3012
$ws.Range("N46").Value = -8155
# Row 61
$ws.Range("H61").Value = 166673000
$ws.Range("I61").Value = 333338000
$ws.Range("J61").Value = 8000
$ws.Range("K61").Value = 333338000
$ws.Range("L61").Value = 8000
$ws.Range("M61").Value = -333337798
$ws.Range("N61").Value = -8404
# Row 68
$ws.Range("H68").Value = 6058.933
$ws.Range("I68").Value = 5372
$ws.Range("J68").Value = 6308.727
$ws.Range("K68").Value = 5372
$ws.Range("L68").Value = 6308.727
$ws.Range("M68").Value = -4623
$ws.Range("N68").Value = -7806.727
# Row 71
$ws.Range("H71").Value = 6058.933
$ws.Range("I71").Value = 5372
$ws.Range("J71").Value = 6308.727
$ws.Range("K71").Value = 26860
$ws.Range("L71").Value = 31543.635
$ws.Range("M71").Value = -23116
$ws.Range("N71").Value = -39031.63499999999
# Row 93
$ws.Range("H93").Value = 4664.3335
$ws.Range("I93").Value = 4622
$ws.Range("J93").Value = 4749
$ws.Range("K93").Value = 4622
$ws.Range("L93").Value = 4749
$ws.Range("M93").Value = -3374
$ws.Range("N93").Value = -7245
# Row 100
$ws.Range("H100").Value = 7607.6665
$ws.Range("I100").Value = 4617.25
$ws.Range("J100").Value = 10000
$ws.Range("K100").Value = 4617.25
$ws.Range("L100").Value = 10000
$ws.Range("M100").Value = -4076.25
$ws.Range("N100").Value = -11082
# Row 113
$ws.Range("H113").Value = 166673000
$ws.Range("I113").Value = 333338000
$ws.Range("J113").Value = 8000
$ws.Range("K113").Value = 333338000
$ws.Range("L113").Value = 8000
$ws.Range("M113").Value = -333335830
$ws.Range("N113").Value = -12340
# Row 127
$ws.Range("H127").Value = 50500
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 50500
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 50500
$ws.Range("N127").Value = -60420

$ws = $wb.Worksheets.Item("WVR")
# Row 95
$ws.Range("H95").Value = 20342.334
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 20342.334
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 20342.334
$ws.Range("N95").Value = -25834.334
# Row 113
$ws.Range("H113").Value = 661.35
$ws.Range("I113").Value = 428.6
$ws.Range("J113").Value = 1359.6
$ws.Range("K113").Value = 1285.8
$ws.Range("L113").Value = 4078.8
$ws.Range("M113").Value = 884.1999999999998
$ws.Range("N113").Value = -8418.799999999999
# Row 126
$ws.Range("H126").Value = 5428.231
$ws.Range("I126").Value = 3118.6
$ws.Range("J126").Value = 6871.75
$ws.Range("K126").Value = 9355.799999999999
$ws.Range("L126").Value = 20615.25
$ws.Range("M126").Value = -6885.799999999999
$ws.Range("N126").Value = -25555.25
# Row 132
$ws.Range("H132").Value = 2995
$ws.Range("I132").Value = 2991.6667
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 8975.000100000001
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -6445.000100000001
$ws.Range("N132").Value = -14060

Write-Host "All changes applied"